$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 header: 0..5 across B8:G8 ---
$headerVals = @(0,1,2,3,4,5)
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item(8, 2 + $i).Value = $headerVals[$i]
}

# --- Column A for rows 9..14: 0..5 ---
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item(9 + $i, 1).Value = $i
}

# --- Formulas (EXP of the original A1:F6 block, shifted) ---
# Typed in B9 first, then filled right across C9:G9 (one shared group),
# then a separate formula typed in B10 and filled down+right across B10:G14
# (second shared group spanning all five remaining rows).
$ws.Range("B9").Formula = "=EXP(A1)"
$ws.Range("C9:G9").Formula = "=EXP(B1)"
$ws.Range("B10:G14").Formula = "=EXP(A2)"

# --- Number format for the new header/index cells: custom "0_);[Red](0)" ---
$customFmt = "0_);[Red]\(0\)"
$ws.Range("B8:G8").NumberFormat = $customFmt
$ws.Range("A9:A14").NumberFormat = $customFmt

# --- Column widths A:C ---
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 9.1640625

# --- Conditional formatting on B9:G14: value >= 0.3 -> red font ---
$cfRange = $ws.Range("B9:G14")
$cf = $cfRange.FormatConditions.Add(1, 7, "0.3")
$cf.Font.Color = 255

# --- Selection ---
$ws.Range("B9").Select()
